$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "FullBlastSMS - Service_List changed to Master_Lists" paragraph:
#    add " Table" (bold, purple 7030A0) right after "Service_List"
#    and " Table" (bold, dark red C00000) right after "Master_Lists".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*FullBlastSMS*Service_List*Master_List*") {

        $pr = $p.Range

        # --- insert " Table" right after "Service_List" ------------------
        $found = $pr.Find.Execute("Service_List", $false, $false, $false, `
                                   $false, $false, $true, 1, $false, "", 0)
        $pr.Collapse(0)                        # wdCollapseEnd
        $pr.MoveEnd(1, 1)                      # grab the following space char
        $pr.Text = " Table "                   # keep it as its own run

        $boldPart = $d.Range($pr.Start, $pr.Start + 6)   # " Table" (6 chars)
        $boldPart.Font.Bold = 1
        $boldPart.Font.Color = 10498160        # RGB(0x70,0x30,0xA0) -> 7030A0

        # --- insert " Table" right after "Master_Lists" / before the
        #     paragraph mark (where the stray _GoBack bookmark used to live)
        $paraTextEnd = $p.Range.End - 1
        $insPoint = $d.Range($paraTextEnd, $paraTextEnd)
        $insPoint.InsertAfter(" Table")
        $insPoint.Font.Bold = 1
        $insPoint.Font.Color = 192             # RGB(0xC0,0x00,0x00) -> C00000

        break
    }
}

# ---------------------------------------------------------------------------
# 2) Drop the _GoBack bookmark from its old spot (right after "Master_Lists").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Add a brand-new (empty) paragraph right after the
#    "Service_Class <tab> Class" paragraph, and park a fresh, collapsed
#    _GoBack bookmark inside it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Service_Class*Class*") {
        $pr = $p.Range
        $textOnlyEnd = $pr.End - 1
        $insPt = $d.Range($textOnlyEnd, $textOnlyEnd)
        $insPt.InsertParagraphAfter()
        break
    }
}

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Service_Class*Class*") {
        $target = $p.Next()
        break
    }
}

$newParaStart = $target.Range.Start
# Placing a bookmark exactly on a bare paragraph mark needs a live
# character to anchor to; insert one, bookmark it, then remove the
# character again so the bookmark collapses in the right spot.
$anchor = $d.Range($newParaStart, $newParaStart)
$anchor.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range($newParaStart, $newParaStart + 1).Text = ""
